$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 112035981
$ws.Range("B2").Value = 90835
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 5964
$ws.Range("F2").Value = "Fjällig taggsvamp s.str."
$ws.Range("G2").Value = "Sarcodon imbricatus s.str."
$ws.Range("H2").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q2").Value = 516149
$ws.Range("R2").Value = 7184413
$ws.Range("S2").Value = 5
$ws.Range("Z2").Value = "11:29"
$ws.Range("AB2").Value = "11:29"
$ws.Range("AH2").Value = "Blåbärsbarrskog"
$ws.Range("AJ2").ClearContents()
$ws.Range("AK2").ClearContents()
$ws.Range("AM2").ClearContents()
$ws.Range("AO2").ClearContents()
$ws.Range("A3").Value = 112035020
$ws.Range("B3").Value = 89549
$ws.Range("E3").Value = 1108
$ws.Range("F3").Value = "Harticka"
$ws.Range("G3").Value = "Pelloporus leporinus"
$ws.Range("H3").Value = "(Fr.) Krieglst."
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("Q3").Value = 515923
$ws.Range("R3").Value = 7184659
$ws.Range("S3").Value = 50
$ws.Range("Z3").Value = "10:24"
$ws.Range("AB3").Value = "10:24"
$ws.Range("AH3").Value = "Blåbärsbarrskog"
$ws.Range("AJ3").ClearContents()
$ws.Range("AK3").ClearContents()
$ws.Range("AM3").Value = "Stubbe"
$ws.Range("AO3").Value = "Stump"
$ws.Range("A4").Value = 112037684
$ws.Range("B4").Value = 77650
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 515886
$ws.Range("R4").Value = 7184226
$ws.Range("Z4").Value = "12:08"
$ws.Range("AB4").Value = "12:08"
$ws.Range("AH4").Value = "Blåbärsbarrskog"
$ws.Range("AM4").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO4").Value = "Standing dead tree/snags # Picea abies"
$ws.Range("A5").Value = 112035549
$ws.Range("B5").Value = 77650
$ws.Range("Q5").Value = 515977
$ws.Range("R5").Value = 7184567
$ws.Range("Z5").Value = "10:51"
$ws.Range("AB5").Value = "10:51"
$ws.Range("AH5").Value = "Blåbärsgranskog"
$ws.Range("AM5").Value = "Gren på levande träd"
$ws.Range("AO5").Value = "Branch on living tree"
$ws.Range("A6").Value = 112037635
$ws.Range("B6").Value = 89549
$ws.Range("E6").Value = 1108
$ws.Range("F6").Value = "Harticka"
$ws.Range("G6").Value = "Pelloporus leporinus"
$ws.Range("H6").Value = "(Fr.) Krieglst."
$ws.Range("Q6").Value = 515886
$ws.Range("R6").Value = 7184226
$ws.Range("Z6").Value = "12:06"
$ws.Range("AB6").Value = "12:06"
$ws.Range("AJ6").Value = "gran"
$ws.Range("AK6").Value = "Picea abies"
$ws.Range("AM6").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO6").Value = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("A7").Value = 112038134
$ws.Range("B7").Value = 89553
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = "Ullticka"
$ws.Range("G7").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H7").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q7").Value = 515925
$ws.Range("R7").Value = 7184319
$ws.Range("S7").Value = 10
$ws.Range("Z7").Value = "13:27"
$ws.Range("AB7").Value = "13:27"
$ws.Range("AH7").Value = "Blåbärsgranskog"
$ws.Range("AJ7").Value = "gran"
$ws.Range("AK7").Value = "Picea abies"
$ws.Range("AM7").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO7").Value = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("A8").Value = 112038529
$ws.Range("B8").Value = 77650
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("Q8").Value = 515872
$ws.Range("R8").Value = 7184628
$ws.Range("Z8").Value = "14:07"
$ws.Range("AB8").Value = "14:07"
$ws.Range("AH8").Value = "Gransumpskog"
$ws.Range("AJ8").ClearContents()
$ws.Range("AK8").ClearContents()
$ws.Range("AM8").ClearContents()
$ws.Range("AO8").ClearContents()
$ws.Range("A9").Value = 112038436
$ws.Range("B9").Value = 89549
$ws.Range("Q9").Value = 515951
$ws.Range("R9").Value = 7184320
$ws.Range("Z9").Value = "13:28"
$ws.Range("AB9").Value = "13:28"
$ws.Range("AM9").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO9").Value = "Standing dead tree/snags # Picea abies"
$ws.Range("A10").Value = 112037386
$ws.Range("B10").Value = 89571
$ws.Range("E10").Value = 5432
$ws.Range("F10").Value = "Granticka"
$ws.Range("G10").Value = "Porodaedalea chrysoloma"
$ws.Range("H10").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q10").Value = 516032
$ws.Range("R10").Value = 7184227
$ws.Range("Z10").Value = "11:52"
$ws.Range("AB10").Value = "11:52"
$ws.Range("AM10").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO10").Value = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("A11").Value = 112037208
$ws.Range("B11").Value = 77650
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("Q11").Value = 516097
$ws.Range("R11").Value = 7184259
$ws.Range("S11").Value = 10
$ws.Range("Z11").Value = "11:44"
$ws.Range("AB11").Value = "11:44"
$ws.Range("AM11").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO11").Value = "Standing dead tree/snags"
$ws.Range("A12").Value = 112038473
$ws.Range("B12").Value = 89834
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 658
$ws.Range("F12").Value = "Rosenticka"
$ws.Range("G12").Value = "Rhodofomes roseus"
$ws.Range("H12").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("I12").Value = "4"
$ws.Range("J12").Value = "fruktkroppar"
$ws.Range("Q12").Value = 516057
$ws.Range("R12").Value = 7184320
$ws.Range("Z12").Value = "13:34"
$ws.Range("AB12").Value = "13:34"
$ws.Range("AM12").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO12").Value = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("A13").Value = 112038082
$ws.Range("B13").Value = 90235
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 3298
$ws.Range("F13").Value = "Trådticka"
$ws.Range("G13").Value = "Climacocystis borealis"
$ws.Range("H13").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q13").Value = 515925
$ws.Range("R13").Value = 7184319
$ws.Range("Z13").Value = "13:22"
$ws.Range("AB13").Value = "13:22"
$ws.Range("AH13").Value = "Blåbärsgranskog"
$ws.Range("AJ13").Value = "gran"
$ws.Range("AK13").Value = "Picea abies"
$ws.Range("AM13").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO13").Value = "Standing dead tree/snags # Picea abies"
